$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'51.441.32"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.30%  "
$ws.Range("D3").Value = "'3.049.05"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'385.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.15%  "
$ws.Range("D6").Value = "'103.15"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.26%  "
$ws.Range("D7").Value = "'0.544"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.35%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "'0.584"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.21%  "
$ws.Range("D10").Value = "'36.76"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.38%  "
$ws.Range("D12").Value = "'0.0861"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.25%  "
$ws.Range("D13").Value = "'3.534.77"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.36%  "
$ws.Range("D14").Value = "'18.54"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.85%  "
$ws.Range("D15").Value = "'7.76"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.38%  "
$ws.Range("D16").Value = "'3.064.47"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.22%  "
$ws.Range("D17").Value = "'0.970"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.85%  "
$ws.Range("D18").Value = "'10.60"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.39%  "
$ws.Range("D19").Value = "'51.541.94"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.35%  "
$ws.Range("D20").Value = "'3.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.43%  "
$ws.Range("D21").Value = "'12.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.95%  "
$ws.Range("E22").Value = "  +0.59%  "
$ws.Range("D23").Value = "'70.12"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.20%  "
$ws.Range("D24").Value = "'267.85"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.26%  "
$ws.Range("D25").Value = "'3.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.71%  "
$ws.Range("E26").Value = "  +5.05%  "
$ws.Range("D27").Value = "'26.90"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.42%  "
$ws.Range("E28").Value = "  +3.20%  "
$ws.Range("D29").Value = "'7.23"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.96%  "
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("E31").Value = "  -1.44%  "
$ws.Range("E32").Value = "  -0.59%  "
$ws.Range("D33").Value = "'34.71"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.24%  "
$ws.Range("E34").Value = "  -0.17%  "
$ws.Range("E35").Value = "  -2.04%  "
$ws.Range("D36").Value = "'0.0448"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.96%  "
$ws.Range("D37").Value = "'1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("D38").Value = "'3.33"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.68%  "
$ws.Range("E39").Value = "  +8.46%  "
$ws.Range("D40").Value = "'16.95"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.41%  "
$ws.Range("D41").Value = "'1.86"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.35%  "
$ws.Range("E42").Value = "  -0.91%  "
$ws.Range("E43").Value = "  +0.14%  "
$ws.Range("D44").Value = "'124.96"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.44%  "
$ws.Range("D45").Value = "'3.75"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.67%  "
$ws.Range("D46").Value = "'21.93"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.35%  "
$ws.Range("D47").Value = "'2.10"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.40%  "
$ws.Range("E48").Value = "  +0.60%  "
$ws.Range("D49").Value = "'2.032.45"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.08%  "
$ws.Range("D50").Value = "'3.350.77"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.49%  "
$ws.Range("B51").Value = "BEAM"
$ws.Range("C51").Value = "https://coinranking.com/coin/cYYMfXF4u+beam-beam"
$ws.Range("D51").Value = "'0.0318"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.50%  "
